$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.765.74'
$ws.Range("E2").Value = '  -0.85%  '

# Row 3
$ws.Range("D3").Value = '2.225.73'
$ws.Range("E3").Value = '  -2.07%  '

# Row 4
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.12'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.88%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.37'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.30%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.567'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.72%  '

# Row 8
$ws.Range("E8").Value = '  +0.15%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -7.55%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.18'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.76%  '

# Row 11
$ws.Range("E11").Value = '  -2.71%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.36'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -6.82%  '

# Row 13
$ws.Range("E13").Value = '  -3.16%  '

# Row 14
$ws.Range("D14").Value = '2.571.27'
$ws.Range("E14").Value = '  -1.96%  '

# Row 15
$ws.Range("D15").Value = '2.228.60'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.835'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.23%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.05'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.00%  '

# Row 18
$ws.Range("D18").Value = '43.673.92'
$ws.Range("E18").Value = '  -0.95%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.99'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -11.03%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0961'
$ws.Range("E20").Value = '  -3.79%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.32'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.40%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.05'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.96%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.98'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -7.39%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '233.37'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.66%  '

# Row 25
$ws.Range("E25").Value = '  -8.96%  '

# Row 26
$ws.Range("E26").Value = '  +0.41%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.12'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.32%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.16'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.88%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.72'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -8.50%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.92'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -9.12%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '157.41'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.31%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.87'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.53%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0829'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -6.86%  '

# Row 34
$ws.Range("E34").Value = '  -1.62%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.19'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.45%  '

# Row 36
$ws.Range("E36").Value = '  +2.24%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.88'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -7.32%  '

# Row 38
$ws.Range("E38").Value = '  -3.72%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.85'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.22%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.59'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -8.87%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.03'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -11.59%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0306'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -6.75%  '

# Row 43
$ws.Range("E43").Value = '  +0.13%  '

# Row 44
$ws.Range("D44").Value = '1.706.30'
$ws.Range("E44").Value = '  -4.90%  '

# Row 45
$ws.Range("E45").Value = '  -7.61%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.11'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.80%  '

# Row 47
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '80.05'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -7.45%  '

# Row 48
$ws.Range("B48").Value = 'ordi'
$ws.Range("C48").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '73.06'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.19%  '

# Row 49
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.76'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.82%  '

# Row 50
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.65'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.14%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.30'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -6.02%  '
